$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: newly committed CRM accuracy data for the sample processed on
# 02/18/2019 (opened bottle from 02/14), taken at Moorea.
$ws.Range("A31").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A32").Value = 43514
$ws.Range("B32").Value = 2216.788
$ws.Range("F32").Value = "New CRM bottle (opened 02/14)"

# Update the active selection to reflect where editing left off.
$ws.Range("G38").Select()
